# Commit: "Added Control System Plot"
# The author added a plot based on the Attributes sheet and, while doing
# so, tweaked two of the lateral/directional derivative inputs (clr, cnb)
# and left the selection on cell R2 (just past the last data column Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Update the two changed attribute values in row 2
$ws.Range("L2").Value = 0.16   # clr
$ws.Range("N2").Value = 0.26   # cnb

# Leave the selection where the author left it after the edit
$ws.Range("R2").Select()
